$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Rename the "Model" row label to "production_function"
$ws.Range("A8").Value = "production_function"

# Insert a new row for the "L_curve" option right below it (default: off)
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0

# Make this sheet the active one (was "threshold_b"), with the selection
# left where editing ended up
$ws.Activate()
$ws.Range("C25").Select()
